$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate formatting of row 32 into the new row 33 (same style/height as existing rows 31/32)
$ws.Range("A32:E32").Copy() | Out-Null
$ws.Range("A33:E33").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Set the new row's content
$ws.Cells.Item(33, 1).Value = "GFG"
$ws.Cells.Item(33, 2).Value = "Longest Sub-Array with Sum K"
$ws.Cells.Item(33, 3).Value = "Java "
$ws.Cells.Item(33, 4).Value = 45003
$ws.Cells.Item(33, 4).NumberFormat = $ws.Cells.Item(32, 4).NumberFormat
$ws.Cells.Item(33, 5).Value = "For +ve and -ve values of element in array and K"

$ws.Rows.Item(33).RowHeight = $ws.Rows.Item(32).RowHeight

# Update the selection to mirror where the user left the cursor after editing
$ws.Range("B36").Select()
